# Updated cryptos list with GitHub Actions
# Applies per-row price (D) and volume-change (E) updates, plus a
# name/link swap between VeChain and EnergySwap (rows 49-50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.021.23"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "2.350.42"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'544.15"
$ws.Range("E5").Value = "  +6.46%  "
$ws.Range("D6").Value = "'134.41"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "2.349.90"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'5.41"
$ws.Range("E12").Value = "  +3.33%  "
$ws.Range("D13").Value = "'0.359"
$ws.Range("E13").Value = "  +6.96%  "
$ws.Range("D14").Value = "2.764.35"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "'23.59"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "58.005.84"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "2.334.41"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("D20").Value = "'333.46"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'61.83"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("D26").Value = "'8.48"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'1.41"
$ws.Range("E28").Value = "  +8.26%  "
$ws.Range("E29").Value = "  +5.34%  "
$ws.Range("D30").Value = "'169.77"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").Value = "'6.14"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("E33").Value = "  +17.50%  "
$ws.Range("D34").Value = "'18.48"
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +7.20%  "
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").Value = "'1.63"
$ws.Range("E39").Value = "  +5.13%  "
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").Value = "'149.03"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.380"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").Value = "'284.30"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("D45").Value = "'19.19"
$ws.Range("E45").Value = "  +6.29%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("D48").Value = "'0.562"
$ws.Range("E48").Value = "  +1.56%  "

# Rows 49 and 50 swap identity (VeChain <-> EnergySwap) and get new
# price / volume-change values.
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'17.61"
$ws.Range("E49").Value = "  +3.92%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0217"
$ws.Range("E50").Value = "  +1.42%  "

# Row 51 price / volume-change update
$ws.Range("D51").Value = "'0.382"
$ws.Range("E51").Value = "  +9.32%  "

Write-Host "cryptos list updated"
